$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 21; $row++) {
    $cell = $ws.Cells.Item($row, 5)  # Column E
    if ($cell.Value2 -eq "fullRNASEQ") {
        $cell.Value = "fullRNASeq"
    }
}
